$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the "items" column (M3) gains an explicit-but-blank text value,
# matching the pattern already used for M2 elsewhere in the sheet. A bare
# apostrophe forces Excel to commit the cell as (blank) text rather than
# leaving it unset.
$ws.Range("M3").Value = "'"

# Row 4: new invoice record appended after row 3.
$ws.Range("A4").Value = "57fe89c5-a399-4dd1-9830-f513fc466f73"
$ws.Range("B4").Value = "INV-0001"

# invoice_date must stay literal text ("2025-10-30"), not get reinterpreted
# as a date serial number.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2025-10-30"

$ws.Range("D4").Value = "draft"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 99
$ws.Range("G4").Value = 8.91
$ws.Range("H4").Value = 8.91
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 116.82
$ws.Range("K4").Value = "chggd"
$ws.Range("L4").Value = "fdhd"
# M4 ("items") intentionally left blank - no line items on this invoice.
$ws.Range("N4").Value = "2025-10-30T17:10:28.233Z"
$ws.Range("O4").Value = "8bf76e3c-b758-48c5-972a-bb86edf041ec"
